# Apply the "Updated symbol list" edit: refresh price/volume(1h) text
# values for each coin row, and swap the BOLO / CoinbaseStockToken rows
# (48 and 49) which changed order in the source feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.24"
$ws.Range("E2").Value = "'1.75%"

$ws.Range("D3").Value = "'41.04"
$ws.Range("E3").Value = "'3.23%"

$ws.Range("D4").Value = "'5.576"
$ws.Range("E4").Value = "'-5.25%"

$ws.Range("D5").Value = "'0.08178"
$ws.Range("E5").Value = "'1.95%"

$ws.Range("D6").Value = "'2.046"
$ws.Range("E6").Value = "'6.37%"

$ws.Range("D7").Value = "'8.733"
$ws.Range("E7").Value = "'0.67%"

$ws.Range("D8").Value = "'4.531"
$ws.Range("E8").Value = "'-0.98%"

$ws.Range("D9").Value = "'2.949"
$ws.Range("E9").Value = "'-0.13%"

$ws.Range("D10").Value = "'0.9182"
$ws.Range("E10").Value = "'-1.33%"

$ws.Range("E11").Value = "'0.78%"

$ws.Range("D12").Value = "'0.1949"
$ws.Range("E12").Value = "'-1.18%"

$ws.Range("D13").Value = "'0.09402"
$ws.Range("E13").Value = "'1.92%"

$ws.Range("D14").Value = "'0.03731"
$ws.Range("E14").Value = "'5.02%"

$ws.Range("E15").Value = "'1.04%"

$ws.Range("D16").Value = "'0.001298"
$ws.Range("E16").Value = "'-0.39%"

$ws.Range("D17").Value = "'0.006212"
$ws.Range("E17").Value = "'0.99%"

$ws.Range("D18").Value = "'3.436"
$ws.Range("E18").Value = "'2.68%"

$ws.Range("E19").Value = "'-2.26%"

$ws.Range("D20").Value = "'8.334"
$ws.Range("E20").Value = "'-4.70%"

$ws.Range("E21").Value = "'-1.80%"

$ws.Range("D22").Value = "'0.2391"
$ws.Range("E22").Value = "'-2.34%"

$ws.Range("D23").Value = "'0.04412"
$ws.Range("E23").Value = "'-0.32%"

$ws.Range("D24").Value = "'0.001263"
$ws.Range("E24").Value = "'0.14%"

$ws.Range("D25").Value = "'0.004304"
$ws.Range("E25").Value = "'-2.35%"

$ws.Range("E26").Value = "'3.68%"

$ws.Range("D39").Value = "'0.02754"
$ws.Range("E39").Value = "'12.48%"

$ws.Range("D40").Value = "'0.05399"

$ws.Range("D41").Value = "'0.007657"
$ws.Range("E41").Value = "'2.67%"

$ws.Range("E42").Value = "'0.64%"

$ws.Range("D43").Value = "'0.009005"
$ws.Range("E43").Value = "'-5.68%"

$ws.Range("D44").Value = "'0.002113"
$ws.Range("E44").Value = "'-0.15%"

$ws.Range("D45").Value = "'0.01126"
$ws.Range("E45").Value = "'13.01%"

$ws.Range("D46").Value = "'0.00006862"
$ws.Range("E46").Value = "'2.08%"

$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.16%"

$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.003581"
$ws.Range("E48").Value = "'19.33%"

$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "'0.002283"
$ws.Range("E49").Value = "'60.46%"

$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.16%"

$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.16%"
